# Auto-generated edit script applying the Durandal_Profits.xlsx data refresh
# (scheduled runner update) across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 17546122
$ws.Range("I100").Value = 41668692
$ws.Range("J100").Value = 2435
$ws.Range("K100").Value = 41668692
$ws.Range("L100").Value = 2435
$ws.Range("M100").Value = -41668151
$ws.Range("N100").Value = -3517

$ws.Range("H120").Value = 37955.4
$ws.Range("J120").Value = 37955.4
$ws.Range("L120").Value = 37955.4
$ws.Range("N120").Value = -47631.4

$ws.Range("H132").Value = 13159019
$ws.Range("I132").Value = 15385161
$ws.Range("J132").Value = 4541.091
$ws.Range("K132").Value = 46155483
$ws.Range("L132").Value = 13623.273
$ws.Range("M132").Value = -46152953
$ws.Range("N132").Value = -18683.273

$ws.Range("H138").Value = 2669.182
$ws.Range("I138").Value = 1676.8605
$ws.Range("J138").Value = 6225
$ws.Range("K138").Value = 5030.5815
$ws.Range("L138").Value = 18675
$ws.Range("M138").Value = 109.4184999999998
$ws.Range("N138").Value = -28955


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1222
$ws.Range("I61").Value = 781.0769
$ws.Range("J61").Value = 1986.2667
$ws.Range("K61").Value = 781.0769
$ws.Range("L61").Value = 1986.2667
$ws.Range("M61").Value = -569.0769
$ws.Range("N61").Value = -2410.2667

$ws.Range("H74").Value = 1073.8611
$ws.Range("I74").Value = 939.8182
$ws.Range("J74").Value = 1284.5
$ws.Range("K74").Value = 939.8182
$ws.Range("L74").Value = 1284.5
$ws.Range("M74").Value = -65.81820000000005
$ws.Range("N74").Value = -3032.5

$ws.Range("H77").Value = 1073.8611
$ws.Range("I77").Value = 939.8182
$ws.Range("J77").Value = 1284.5
$ws.Range("K77").Value = 4699.091
$ws.Range("L77").Value = 6422.5
$ws.Range("M77").Value = -331.0910000000003
$ws.Range("N77").Value = -15158.5

$ws.Range("H132").Value = 2447
$ws.Range("I132").Value = 1271
$ws.Range("J132").Value = 3740.6
$ws.Range("K132").Value = 3813
$ws.Range("L132").Value = 11221.8
$ws.Range("M132").Value = -1283
$ws.Range("N132").Value = -16281.8

$ws.Range("H136").Value = 1222
$ws.Range("I136").Value = 781.0769
$ws.Range("J136").Value = 1986.2667
$ws.Range("K136").Value = 2343.2307
$ws.Range("L136").Value = 5958.800099999999
$ws.Range("M136").Value = 206.7692999999999
$ws.Range("N136").Value = -11058.8001


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4484.05
$ws.Range("I134").Value = 720.65625
$ws.Range("J134").Value = 19537.625
$ws.Range("K134").Value = 2161.96875
$ws.Range("L134").Value = 58612.875
$ws.Range("M134").Value = 373.03125
$ws.Range("N134").Value = -63682.875


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5815499.5
$ws.Range("I31").Value = 8334601
$ws.Range("J31").Value = 2188.8845
$ws.Range("K31").Value = 8334601
$ws.Range("L31").Value = 2188.8845
$ws.Range("M31").Value = -8334306
$ws.Range("N31").Value = -2778.8845

$ws.Range("H34").Value = 5815499.5
$ws.Range("I34").Value = 8334601
$ws.Range("J34").Value = 2188.8845
$ws.Range("K34").Value = 8334601
$ws.Range("L34").Value = 2188.8845
$ws.Range("M34").Value = -8334399
$ws.Range("N34").Value = -2592.8845

$ws.Range("H58").Value = 1035.8788
$ws.Range("I58").Value = 741.2857
$ws.Range("J58").Value = 2685.6
$ws.Range("K58").Value = 741.2857
$ws.Range("L58").Value = 2685.6
$ws.Range("M58").Value = -538.2857
$ws.Range("N58").Value = -3091.6

$ws.Range("H99").Value = 1319.2413
$ws.Range("I99").Value = 1310.4
$ws.Range("J99").Value = 1374.5
$ws.Range("K99").Value = 1310.4
$ws.Range("L99").Value = 1374.5
$ws.Range("M99").Value = 187.5999999999999
$ws.Range("N99").Value = -4370.5

$ws.Range("H126").Value = 1319.2413
$ws.Range("I126").Value = 1310.4
$ws.Range("J126").Value = 1374.5
$ws.Range("K126").Value = 3931.2
$ws.Range("L126").Value = 4123.5
$ws.Range("M126").Value = -1461.2
$ws.Range("N126").Value = -9063.5

$ws.Range("H132").Value = 34736.9
$ws.Range("I132").Value = 883.8
$ws.Range("J132").Value = 204002.4
$ws.Range("K132").Value = 2651.4
$ws.Range("L132").Value = 612007.2
$ws.Range("M132").Value = -121.3999999999996
$ws.Range("N132").Value = -617067.2

$ws.Range("H133").Value = 100000
$ws.Range("J133").Value = 100000
$ws.Range("L133").Value = 100000
$ws.Range("N133").Value = -105060

$ws.Range("H135").Value = 42857.145
$ws.Range("J135").Value = 42857.145
$ws.Range("L135").Value = 42857.145
$ws.Range("N135").Value = -52997.145

$ws.Range("H136").Value = 1035.8788
$ws.Range("I136").Value = 741.2857
$ws.Range("J136").Value = 2685.6
$ws.Range("K136").Value = 2223.8571
$ws.Range("L136").Value = 8056.799999999999
$ws.Range("M136").Value = 326.1428999999998
$ws.Range("N136").Value = -13156.8

$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws.Range("H141").Value = 22640
$ws.Range("J141").Value = 22640
$ws.Range("L141").Value = 22640
$ws.Range("N141").Value = -33000


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 447.03845
$ws.Range("I107").Value = 322.83334
$ws.Range("J107").Value = 553.5
$ws.Range("K107").Value = 968.5000200000001
$ws.Range("L107").Value = 1660.5
$ws.Range("M107").Value = 951.4999799999999
$ws.Range("N107").Value = -5500.5

$ws.Range("H119").Value = 7330.9287
$ws.Range("I119").Value = 6401.5
$ws.Range("J119").Value = 8028
$ws.Range("K119").Value = 19204.5
$ws.Range("L119").Value = 24084
$ws.Range("M119").Value = -14366.5
$ws.Range("N119").Value = -33760

$ws.Range("H132").Value = 2496.1155
$ws.Range("I132").Value = 1542.8572
$ws.Range("J132").Value = 2847.3157
$ws.Range("K132").Value = 13885.7148
$ws.Range("L132").Value = 25625.8413
$ws.Range("M132").Value = -11355.7148
$ws.Range("N132").Value = -30685.8413


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 113901.89
$ws.Range("I132").Value = 102532.2
$ws.Range("J132").Value = 128114
$ws.Range("K132").Value = 307596.6
$ws.Range("L132").Value = 384342
$ws.Range("M132").Value = -305066.6
$ws.Range("N132").Value = -389402


# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2228.3442
$ws.Range("I132").Value = 2169.4314
$ws.Range("J132").Value = 2528.8
$ws.Range("K132").Value = 6508.2942
$ws.Range("L132").Value = 7586.400000000001
$ws.Range("M132").Value = -3978.2942
$ws.Range("N132").Value = -12646.4

$ws.Range("H136").Value = 2378
$ws.Range("I136").Value = 1270.0625
$ws.Range("J136").Value = 7442.857
$ws.Range("K136").Value = 3810.1875
$ws.Range("L136").Value = 22328.571
$ws.Range("M136").Value = -1260.1875
$ws.Range("N136").Value = -27428.571


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 1302
$ws.Range("J5").Value = 1302
$ws.Range("L5").Value = 1302
$ws.Range("N5").Value = -1526

$ws.Range("H132").Value = 19534218
$ws.Range("I132").Value = 27778966
$ws.Range("J132").Value = 7183.316
$ws.Range("K132").Value = 83336898
$ws.Range("L132").Value = 21549.948
$ws.Range("M132").Value = -83334368
$ws.Range("N132").Value = -26609.948

